$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 08:32:32"
$ws1.Range("A3").Value = "Total filas: 122"

$ws1data = @(
  @("00:09:57", "01:11", "215_ALUAR", 62, "LP1912"),
  @("01:17:09", "01:58", "14_ABASTO", 41, "LP1912"),
  @("01:59:40", "02:00", "14_ABASTO", 1, "LP1912"),
  @("02:30:53", "02:57", "215_ALUAR", 27, "LP1912"),
  @("02:55:01", "02:58", "215_ALUAR", 3, "LP1912"),
  @("01:59:40", "03:06", "215_ALUAR", 67, "LP1912"),
  @("03:21:41", "03:24", "14_ABASTO", 3, "LP1912"),
  @("02:55:01", "03:48", "14_ABASTO", 53, "LP1912"),
  @("01:59:40", "03:50", "14_ABASTO", 111, "LP1912"),
  @("02:30:53", "03:52", "14_ABASTO", 82, "LP1912"),
  @("03:57:17", "04:01", "81_EL PELIGRO", 4, "LP1912"),
  @("03:21:41", "04:45", "215A_EL PATO", 84, "LP1912"),
  @("04:46:27", "04:46", "215A_EL PATO", 0, "LP1912"),
  @("04:46:27", "04:53", "11_ETCHEVERRY", 7, "LP1912"),
  @("03:21:41", "05:14", "14_ABASTO", 113, "LP1912"),
  @("04:56:17", "05:16", "17_ROMERO", 20, "LP1912"),
  @("04:56:17", "05:22", "23_HERNANDEZ", 26, "LP1912"),
  @("05:23:05", "05:24", "23_HERNANDEZ", 1, "LP1912"),
  @("04:46:27", "05:34", "215B_EL PATO", 48, "LP1912"),
  @("05:23:05", "05:35", "215B_EL PATO", 12, "LP1912"),
  @("03:57:17", "05:39", "14_ABASTO", 102, "LP1912"),
  @("05:23:05", "05:46", "15_ABASTO", 23, "LP1912"),
  @("04:24:09", "05:48", "215A_EL PATO", 84, "LP1912"),
  @("05:54:50", "05:54", "10_OLMOS", 0, "LP1912"),
  @("05:54:50", "06:04", "16_SANTA ANA", 10, "LP1912"),
  @("04:24:09", "06:09", "16_SANTA ANA", 105, "LP1912"),
  @("05:54:50", "06:11", "215A_EL PATO", 17, "LP1912"),
  @("05:23:05", "06:12", "215A_EL PATO", 49, "LP1912"),
  @("05:54:50", "06:14", "225_HARAS DEL SUR", 20, "LP1912"),
  @("05:54:50", "06:21", "26_HERNANDEZ", 27, "LP1912"),
  @("05:54:50", "06:27", "23_HERNANDEZ", 33, "LP1912"),
  @("06:24:49", "06:29", "86_EST CHICA-ESC AGRARIA", 5, "LP1912"),
  @("05:23:05", "06:30", "86_EST CHICA-ESC AGRARIA", 67, "LP1912"),
  @("06:24:49", "06:31", "16_SANTA ANA", 7, "LP1912"),
  @("06:24:49", "06:32", "23_HERNANDEZ", 8, "LP1912"),
  @("06:24:49", "06:44", "225_C ROCA-H SUR", 20, "LP1912"),
  @("06:24:49", "06:46", "215C_EL PATO", 22, "LP1912"),
  @("05:23:05", "06:47", "215C_EL PATO", 84, "LP1912"),
  @("06:52:34", "06:59", "14_ABASTO", 7, "LP1912"),
  @("05:23:05", "07:00", "14_ABASTO", 97, "LP1912"),
  @("06:24:49", "07:01", "16_SANTA ANA", 37, "LP1912"),
  @("05:54:50", "07:04", "23_HERNANDEZ", 70, "LP1912"),
  @("06:52:34", "07:05", "15_ABASTO", 13, "LP1912"),
  @("06:52:34", "07:05", "23_HERNANDEZ", 13, "LP1912"),
  @("06:24:49", "07:06", "225_GOMEZ", 42, "LP1912"),
  @("06:52:34", "07:07", "225_GOMEZ", 15, "LP1912"),
  @("06:52:34", "07:11", "215A_EL PATO", 19, "LP1912"),
  @("05:23:05", "07:12", "215A_EL PATO", 109, "LP1912"),
  @("06:52:34", "07:15", "11_ETCHEVERRY", 23, "LP1912"),
  @("06:52:34", "07:16", "16_SANTA ANA", 24, "LP1912"),
  @("05:23:05", "07:16", "11_ETCHEVERRY", 113, "LP1912"),
  @("07:18:13", "07:20", "10_OLMOS", 2, "LP1912"),
  @("07:18:13", "07:21", "26_HERNANDEZ", 3, "LP1912"),
  @("06:52:34", "07:23", "10_OLMOS", 31, "LP1912"),
  @("06:52:34", "07:31", "11_ETCHEVERRY", 39, "LP1912"),
  @("07:18:13", "07:32", "84_COLONIA URQUIZA-ESC 49", 14, "LP1912"),
  @("07:18:13", "07:32", "11_ETCHEVERRY", 14, "LP1912"),
  @("07:18:13", "07:35", "23_HERNANDEZ", 17, "LP1912"),
  @("06:52:34", "07:36", "27_EL RETIRO", 44, "LP1912"),
  @("07:18:13", "07:37", "27_EL RETIRO", 19, "LP1912"),
  @("07:18:13", "07:39", "10_OLMOS", 21, "LP1912"),
  @("05:54:50", "07:46", "16_SANTA ANA", 112, "LP1912"),
  @("06:52:34", "07:47", "14_ABASTO", 55, "LP1912"),
  @("07:18:13", "07:48", "14_ABASTO", 30, "LP1912"),
  @("07:49:14", "07:49", "10_OLMOS", 0, "LP1912"),
  @("07:49:14", "07:51", "215D_EL PATO", 2, "LP1912"),
  @("07:18:13", "07:52", "215D_EL PATO", 34, "LP1912"),
  @("06:52:34", "07:58", "16_SANTA ANA", 66, "LP1912"),
  @("07:49:14", "07:59", "16_SANTA ANA", 10, "LP1912"),
  @("07:49:14", "08:00", "23_HERNANDEZ", 11, "LP1912"),
  @("06:24:49", "08:05", "23_HERNANDEZ", 101, "LP1912"),
  @("07:18:13", "08:05", "11_ETCHEVERRY", 47, "LP1912"),
  @("06:52:34", "08:06", "23_HERNANDEZ", 74, "LP1912"),
  @("07:49:14", "08:08", "11_ETCHEVERRY", 19, "LP1912"),
  @("08:02:29", "08:10", "11_ETCHEVERRY", 8, "LP1912"),
  @("07:49:14", "08:11", "16_SANTA ANA", 22, "LP1912"),
  @("08:02:29", "08:12", "15_ABASTO", 10, "LP1912"),
  @("08:02:29", "08:13", "10_OLMOS", 11, "LP1912"),
  @("08:02:29", "08:14", "16_SANTA ANA", 12, "LP1912"),
  @("08:02:29", "08:21", "26_HERNANDEZ", 19, "LP1912"),
  @("08:02:29", "08:22", "16_P MOR-SANTA ANA", 20, "LP1912"),
  @("07:18:13", "08:23", "16_SANTA ANA", 65, "LP1912"),
  @("07:49:14", "08:23", "16_P MOR-SANTA ANA", 34, "LP1912"),
  @("08:02:29", "08:23", "215B_EL PATO", 21, "LP1912"),
  @("08:02:29", "08:27", "84_COLONIA URQUIZA-ESC 49", 25, "LP1912"),
  @("07:49:14", "08:31", "23_HERNANDEZ", 42, "LP1912"),
  @("08:02:29", "08:33", "10_OLMOS", 31, "LP1912"),
  @("08:02:29", "08:34", "23_HERNANDEZ", 32, "LP1912"),
  @("08:32:32", "08:41", "23_HERNANDEZ", 9, "LP1912"),
  @("08:32:32", "08:42", "81_EL PELIGRO", 10, "LP1912"),
  @("08:32:32", "08:43", "14_ABASTO", 11, "LP1912"),
  @("07:49:14", "08:44", "14_ABASTO", 55, "LP1912"),
  @("08:32:32", "08:53", "10_OLMOS", 21, "LP1912"),
  @("08:32:32", "08:54", "17_ROMERO", 22, "LP1912"),
  @("08:32:32", "09:01", "215A_EL PATO", 29, "LP1912"),
  @("07:49:14", "09:02", "215A_EL PATO", 73, "LP1912"),
  @("08:32:32", "09:03", "11_ETCHEVERRY", 31, "LP1912"),
  @("08:32:32", "09:04", "23_HERNANDEZ", 32, "LP1912"),
  @("07:49:14", "09:04", "11_ETCHEVERRY", 75, "LP1912"),
  @("08:02:29", "09:08", "23_HERNANDEZ", 66, "LP1912"),
  @("08:32:32", "09:09", "16_SANTA ANA", 37, "LP1912"),
  @("08:32:32", "09:10", "16_P MOR-SANTA ANA", 38, "LP1912"),
  @("07:49:14", "09:11", "16_P MOR-SANTA ANA", 82, "LP1912"),
  @("08:32:32", "09:16", "27_EL RETIRO", 44, "LP1912"),
  @("07:49:14", "09:17", "27_EL RETIRO", 88, "LP1912"),
  @("08:02:29", "09:21", "16_SANTA ANA", 79, "LP1912"),
  @("08:32:32", "09:21", "26_HERNANDEZ", 49, "LP1912"),
  @("08:32:32", "09:22", "16_SANTA ANA", 50, "LP1912"),
  @("08:32:32", "09:22", "17_ROMERO", 50, "LP1912"),
  @("08:02:29", "09:23", "17_ROMERO", 81, "LP1912"),
  @("08:32:32", "09:23", "11_ETCHEVERRY", 51, "LP1912"),
  @("07:49:14", "09:24", "11_ETCHEVERRY", 95, "LP1912"),
  @("08:32:32", "09:32", "15_ABASTO", 60, "LP1912"),
  @("08:32:32", "09:33", "10_OLMOS", 61, "LP1912"),
  @("08:32:32", "09:42", "215C_EL PATO", 70, "LP1912"),
  @("08:32:32", "09:43", "14_ABASTO", 71, "LP1912"),
  @("07:49:14", "09:44", "14_ABASTO", 115, "LP1912"),
  @("08:32:32", "09:47", "10_OLMOS", 75, "LP1912"),
  @("08:32:32", "10:10", "16_P MOR-SANTA ANA", 98, "LP1912"),
  @("08:32:32", "10:12", "15_ABASTO", 100, "LP1912"),
  @("08:32:32", "10:21", "26_HERNANDEZ", 109, "LP1912"),
  @("08:32:32", "10:26", "215A_EL PATO", 114, "LP1912")
)

for ($i = 0; $i -lt $ws1data.Count; $i++) {
    $row = $ws1data[$i]
    $r = 6 + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws1.Cells.Item($r, $j+1).Value = $row[$j]
    }
}

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:32:32"
$ws2.Range("A3").Value = "Total filas: 22"

$ws2data = @(
  @("00:09:57", "01:11", "215_ALUAR", 62, "LP1912"),
  @("02:30:53", "02:57", "215_ALUAR", 27, "LP1912"),
  @("02:55:01", "02:58", "215_ALUAR", 3, "LP1912"),
  @("01:59:40", "03:06", "215_ALUAR", 67, "LP1912"),
  @("03:21:41", "04:45", "215A_EL PATO", 84, "LP1912"),
  @("04:46:27", "04:46", "215A_EL PATO", 0, "LP1912"),
  @("04:46:27", "05:34", "215B_EL PATO", 48, "LP1912"),
  @("05:23:05", "05:35", "215B_EL PATO", 12, "LP1912"),
  @("04:24:09", "05:48", "215A_EL PATO", 84, "LP1912"),
  @("05:54:50", "06:11", "215A_EL PATO", 17, "LP1912"),
  @("05:23:05", "06:12", "215A_EL PATO", 49, "LP1912"),
  @("06:24:49", "06:46", "215C_EL PATO", 22, "LP1912"),
  @("05:23:05", "06:47", "215C_EL PATO", 84, "LP1912"),
  @("06:52:34", "07:11", "215A_EL PATO", 19, "LP1912"),
  @("05:23:05", "07:12", "215A_EL PATO", 109, "LP1912"),
  @("07:49:14", "07:51", "215D_EL PATO", 2, "LP1912"),
  @("07:18:13", "07:52", "215D_EL PATO", 34, "LP1912"),
  @("08:02:29", "08:23", "215B_EL PATO", 21, "LP1912"),
  @("08:32:32", "09:01", "215A_EL PATO", 29, "LP1912"),
  @("07:49:14", "09:02", "215A_EL PATO", 73, "LP1912"),
  @("08:32:32", "09:42", "215C_EL PATO", 70, "LP1912"),
  @("08:32:32", "10:26", "215A_EL PATO", 114, "LP1912")
)

for ($i = 0; $i -lt $ws2data.Count; $i++) {
    $row = $ws2data[$i]
    $r = 6 + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws2.Cells.Item($r, $j+1).Value = $row[$j]
    }
}

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:32:32"
$ws3.Range("A3").Value = "Total filas: 19"

$ws3data = @(
  @("00:09:57", "00:09", "215A_LA PLATA", 0, "L6173"),
  @("05:23:05", "05:44", "215A_LA PLATA", 21, "L6173"),
  @("05:54:50", "06:08", "215A_LA PLATA", 14, "L6173"),
  @("05:23:05", "06:09", "215A_LA PLATA", 46, "L6173"),
  @("06:24:49", "06:32", "215C_LA PLATA", 8, "L6203"),
  @("05:23:05", "06:33", "215C_LA PLATA", 70, "L6203"),
  @("06:24:49", "06:59", "215B_LP-P MOR-1 Y 57", 35, "L6173"),
  @("06:52:34", "07:00", "215B_LP-P MOR-1 Y 57", 8, "L6173"),
  @("07:18:13", "07:35", "215A_LA PLATA", 17, "L6173"),
  @("06:52:34", "08:06", "215C_LA PLATA", 74, "L6203"),
  @("07:18:13", "08:08", "215C_LA PLATA", 50, "L6203"),
  @("07:49:14", "08:11", "215C_LA PLATA", 22, "L6203"),
  @("08:02:29", "08:13", "215C_LA PLATA", 11, "L6203"),
  @("06:52:34", "08:33", "215A_LA PLATA", 101, "L6173"),
  @("07:49:14", "08:35", "215A_LA PLATA", 46, "L6173"),
  @("08:02:29", "08:37", "215A_LA PLATA", 35, "L6173"),
  @("08:32:32", "08:38", "215A_LA PLATA", 6, "L6173"),
  @("08:32:32", "09:09", "215D_LA PLATA", 37, "L6203"),
  @("08:32:32", "10:03", "215B_LP-P MOR-40 Y 115", 91, "L6173")
)

for ($i = 0; $i -lt $ws3data.Count; $i++) {
    $row = $ws3data[$i]
    $r = 6 + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws3.Cells.Item($r, $j+1).Value = $row[$j]
    }
}
